$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.052.73'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').Value = '2.668.28'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('D5').Value = '524.15'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '144.48'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '0.569'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Value = '6.97'
$ws.Range('E9').Value = '  +7.71%  '
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('D11').Value = '0.336'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = '3.133.35'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').Value = '59.046.63'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '21.05'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '2.677.85'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').Value = "'339.50"
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').Value = "'4.40"
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('D20').Value = '10.38'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '64.45'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').Value = "'0.420"
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('D27').Value = '0.0₃0803'
$ws.Range('D28').Value = '7.15'
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('D29').Value = '6.68'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = "'1.60"
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').Value = '18.92'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').Value = '150.64'
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('D34').Value = '4.16'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('D36').Value = '0.902'
$ws.Range('E36').Value = '  -5.15%  '
$ws.Range('D37').Value = '0.873'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  -5.90%  '
$ws.Range('D40').Value = '3.59'
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('D41').Value = '0.618'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').Value = '275.62'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').Value = '19.84'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('E45').Value = '  -1.83%  '
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').Value = '2.054.29'
$ws.Range('E48').Value = '  -3.78%  '
$ws.Range('D49').Value = '4.73'
$ws.Range('E49').Value = '  -3.19%  '
$ws.Range('E50').Value = '  -2.78%  '
$ws.Range('D51').Value = '18.89'
$ws.Range('E51').Value = '  -1.75%  '
